$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Marking): Right count 4 -> 5, Wrong marking -1 -> -1.2
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 (Total): Right score 56 -> 70, Wrong marking -1 -> -1.2, Max text 55/112 -> 68.8/140
$ws.Range("B12").Value = 70
$ws.Range("C12").Value = -1.2
$ws.Range("E12").Value = "68.8/140"
